# Add the new CRM accuracy data row (run 7, CRM opened 2021-06-24),
# extending the existing logged titration data by one entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A43").Value = 20210624
$ws.Range("B43").Value = 2228.4520000000002
$ws.Range("C43").Value = 2224.4699999999998
$ws.Range("D43").Formula = "=100*(B43-C43)/C43"
$ws.Range("E43").Value = 180
$ws.Range("F43").Value = "CRM OPENED 20210624"

# Match the author's final scroll/cursor position from the commit.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 1
$ws.Range("H44").Select()
